$d = $word.ActiveDocument

# 1. Title: "Conditional Probability... and Magic!" -> French
$r = $d.Content
$r.Find.Execute("Conditional Probability... and Magic!")
$r.Text = "Probabilité conditionnelle... et magie !"

# 2. "Mathematical Thinking" -> "Pensée mathématique"
$r = $d.Content
$r.Find.Execute("Mathematical Thinking")
$r.Text = "Pensée mathématique"

# 3. Long description paragraph
$r = $d.Content
$r.Find.Execute("What you know about a problem guides you towards the correct solution: conditional probability tells you why. But there is more: let’s discover together how the same concept can be applied to less intuitive – and somewhat magical - contexts.")
$r.Text = "Ce que vous savez d'un problème vous oriente vers la bonne solution : la probabilité conditionnelle vous explique pourquoi. Mais ce n'est pas tout : découvrons ensemble comment le même concept peut être appliqué à des contextes moins intuitifs - et quelque peu magiques."

# 4. "45min" -> "45 minutes"
$r = $d.Content
$r.Find.Execute("45min")
$r.Text = "45 minutes"

# 5. "N. des étudiants" -> "N. Des étudiants" (capitalisation fix)
$r = $d.Content
$r.Find.Execute("N. des étudiants")
$r.Text = "N. Des étudiants"

# 6. Materials sentence (partial translation, "A board would be helpful." stays in English)
$r = $d.Content
$r.Find.Execute("Per group of learners: deck of cards, 2-3 blank paper sheets, and a pen. A board would be helpful.")
$r.Text = "Par groupe d'apprenants : jeu de cartes, 2 ou 3 feuilles de papier vierge et un stylo. A board would be helpful."
